# Apply the "GW_10 / Content info help-link" test-case addition to the check list sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 10 (existing "GW_9" row): change Sub-Module from "Main Search" to
# "Content info", update the "Expected Results" text (F10) and highlight the
# "No." cell (A10) with the same green fill used for GW_8/GW_9 (A8/A9).
# ---------------------------------------------------------------------------
$ws.Range("D10").Value = "Content info"

$ws.Range("F10").Value = "`n1. Website correctly open on Chrome browser.`n2. The entered text is displayed correctly in the input field.`n3. Google search page with query results has special panel 'content info' - footer for page"

# Match A8/A9 styling (green fill FF92D050, same font/alignment as before).
$ws.Range("A10").Interior.Color = 5296274

# ---------------------------------------------------------------------------
# Row 11 (previously a blank filler row): new "GW_10" test case about the
# working 'help' link in the content info panel.
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "GW_10"
$ws.Range("C11").Value = "compatibility"
$ws.Range("D11").Value = "Content info"

# Match the "compatibility" column styling used by the rest of the table
# (centered, bold Arial 10pt) by copying the format from the cell right
# above it.
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("E11").Value = "Check that the panel 'content info' has working link 'help' after opening the serch results page:`n1. Open 'https://www.google.com/' with Chrome`n2. Enter a valid query in the search (for example: wikipedia). `n3. Use key 'Enter'. `n4. Click on the link 'help' at the bottom of the screen"
$boldLen = "Check that the panel 'content info' has working link 'help' after opening the serch results page:".Length
$ws.Range("E11").Characters(1, $boldLen).Font.Bold = $true

$ws.Range("F11").Value = "`n1. Website correctly open on Chrome browser.`n2. The entered text is displayed correctly in the input field.`n3. Google search page with query results correctly open.`n4. Google search help correctly opens"

$ws.Rows.Item(11).RowHeight = 155.25

# ---------------------------------------------------------------------------
# Misc: keep the active-cell selection in sync with the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("G13").Select() | Out-Null
